# Spelling correction in Tokyo
# Tokyo, Hyogo, Osaka and Kyoto were corrected in the dictionary to remove the long o's.
# As a result, the Visit_hx (column F) values that were previously "Unknown" for the
# affected prefectures are now resolved to their correct Yes/No values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F14").Value = "No"    # Hyogo
$ws.Range("F21").Value = "No"    # Kochi
$ws.Range("F23").Value = "Yes"   # Kyoto
$ws.Range("F31").Value = "No"    # Oita
$ws.Range("F34").Value = "Yes"   # Osaka
$ws.Range("F42").Value = "Yes"   # Tokyo
